# Append a new row (row 10) to Sheet1, recording another trip/ration
# entry, same shape as the existing rows in the log.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("ملاحظات" / notes) is left blank for this entry, just like
# most of the existing rows (so we don't touch Range("A10") at all).
$ws.Range("B10").Value = "احمد"

# Column C ("الكمية" / quantity) holds a number-looking value but the
# source data stores it as text, so force a Text number format before
# assigning the value - otherwise Excel would coerce "2" into a real
# number.
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "2"

$ws.Range("D10").Value = "الجزائري"
$ws.Range("E10").Value = "الرحلة 1"
$ws.Range("F10").Value = "C3"
$ws.Range("G10").Value = "NRC"
$ws.Range("H10").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٤٢:٠٤ م"
